# Refresh market-price-derived columns (H-N) across all Leve profit sheets.
# Values below mirror the latest Universalis price snapshot pulled by the
# scheduled runner; only numeric market columns are touched, nothing else.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41: The Write Stuff
$ws.Range("H41").Value = 1310.619
$ws.Range("J41").Value = 1035.909
$ws.Range("L41").Value = 1035.909
$ws.Range("N41").Value = -1915.909

# Row 53: No Accounting for Waste
$ws.Range("H53").Value = 100
$ws.Range("I53").Value = 100
$ws.Range("K53").Value = 100
$ws.Range("M53").Value = 537

# Row 69: Steeling the Knife, Steeling the Mind
$ws.Range("H69").Value = 17999.8
$ws.Range("I69").Value = 23333.334
$ws.Range("J69").Value = 15714
$ws.Range("K69").Value = 70000.00199999999
$ws.Range("L69").Value = 47142
$ws.Range("M69").Value = -69126.00199999999
$ws.Range("N69").Value = -48890

# Row 72: Surgical Substitution (L)
$ws.Range("H72").Value = 17999.8
$ws.Range("I72").Value = 23333.334
$ws.Range("J72").Value = 15714
$ws.Range("K72").Value = 210000.006
$ws.Range("L72").Value = 141426
$ws.Range("M72").Value = -205632.006
$ws.Range("N72").Value = -150162

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 88893384
$ws.Range("I86").Value = 66668668
$ws.Range("K86").Value = 66668668
$ws.Range("M86").Value = -66667545

# Row 87: There Was a Late Fee
$ws.Range("H87").Value = 64333.168
$ws.Range("J87").Value = 64333.168
$ws.Range("L87").Value = 64333.168
$ws.Range("N87").Value = -66829.16800000001

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 88893384
$ws.Range("I89").Value = 66668668
$ws.Range("K89").Value = 333343340
$ws.Range("M89").Value = -333337724

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("H90").Value = 64333.168
$ws.Range("J90").Value = 64333.168
$ws.Range("L90").Value = 192999.504
$ws.Range("N90").Value = -205479.504

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 566.1429000000001
$ws.Range("I118").Value = 595.9167
$ws.Range("K118").Value = 1787.7501
$ws.Range("M118").Value = -130.7501

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6716.0923
$ws.Range("J32").Value = 19706.625
$ws.Range("L32").Value = 19706.625
$ws.Range("N32").Value = -20280.625

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2265.2307
$ws.Range("I61").Value = 1507.3684
$ws.Range("K61").Value = 1507.3684
$ws.Range("M61").Value = -1295.3684

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2293.1396
$ws.Range("I74").Value = 1451.9
$ws.Range("J74").Value = 3024.652
$ws.Range("K74").Value = 1451.9
$ws.Range("L74").Value = 3024.652
$ws.Range("M74").Value = -577.9000000000001
$ws.Range("N74").Value = -4772.652

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2293.1396
$ws.Range("I77").Value = 1451.9
$ws.Range("J77").Value = 3024.652
$ws.Range("K77").Value = 7259.5
$ws.Range("L77").Value = 15123.26
$ws.Range("M77").Value = -2891.5
$ws.Range("N77").Value = -23859.26

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 643.6923
$ws.Range("I110").Value = 643.6923
$ws.Range("K110").Value = 643.6923
$ws.Range("M110").Value = 1401.3077

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1929.8485
$ws.Range("I132").Value = 1484.6296
$ws.Range("K132").Value = 4453.8888
$ws.Range("M132").Value = -1923.8888

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2265.2307
$ws.Range("I136").Value = 1507.3684
$ws.Range("K136").Value = 4522.1052
$ws.Range("M136").Value = -1972.1052

$ws = $wb.Worksheets.Item("BSM")
# Row 81: Diamond Sawdust
$ws.Range("H81").Value = 19855.285
$ws.Range("J81").Value = 19831.166
$ws.Range("L81").Value = 19831.166
$ws.Range("N81").Value = -21953.166

# Row 84: I'm a Lumberjack and I'm Okay (L)
$ws.Range("H84").Value = 19855.285
$ws.Range("J84").Value = 19831.166
$ws.Range("L84").Value = 59493.49800000001
$ws.Range("N84").Value = -70101.49800000001

# Row 96: Hammer Time
$ws.Range("H96").Value = 10988
$ws.Range("I96").Value = 10988
$ws.Range("K96").Value = 10988
$ws.Range("M96").Value = -8242

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 57634.223
$ws.Range("J105").Value = 4183.1665
$ws.Range("L105").Value = 4183.1665
$ws.Range("N105").Value = -7677.1665

# Row 107: The Gold Experience
$ws.Range("H107").Value = 2044.7059
$ws.Range("I107").Value = 1504.2222
$ws.Range("K107").Value = 1504.2222
$ws.Range("M107").Value = 415.7778000000001

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2086.6843
$ws.Range("I134").Value = 1217.2413
$ws.Range("K134").Value = 3651.7239
$ws.Range("M134").Value = -1116.7239

$ws = $wb.Worksheets.Item("CRP")
# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 44427.375
$ws.Range("I51").Value = 38884
$ws.Range("J51").Value = 53666.332
$ws.Range("K51").Value = 38884
$ws.Range("L51").Value = 53666.332
$ws.Range("M51").Value = -38148
$ws.Range("N51").Value = -55138.332

# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 26187.375
$ws.Range("I60").Value = 11624.75
$ws.Range("J60").Value = 40750
$ws.Range("K60").Value = 11624.75
$ws.Range("L60").Value = 40750
$ws.Range("M60").Value = -11113.75
$ws.Range("N60").Value = -41772

# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 44427.375
$ws.Range("I61").Value = 38884
$ws.Range("J61").Value = 53666.332
$ws.Range("K61").Value = 38884
$ws.Range("L61").Value = 53666.332
$ws.Range("M61").Value = -38536
$ws.Range("N61").Value = -54362.332

# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 52500
$ws.Range("J68").Value = 57500
$ws.Range("L68").Value = 57500
$ws.Range("N68").Value = -58998

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 52500
$ws.Range("J71").Value = 57500
$ws.Range("L71").Value = 172500
$ws.Range("N71").Value = -179988

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 1342.3158
$ws.Range("J5").Value = 2003.1428
$ws.Range("L5").Value = 6009.428400000001
$ws.Range("N5").Value = -6233.428400000001

# Row 132: More Mezcal
$ws.Range("H132").Value = 3279.4666
$ws.Range("I132").Value = 2350.5
$ws.Range("J132").Value = 3898.7778
$ws.Range("K132").Value = 21154.5
$ws.Range("L132").Value = 35089.00019999999
$ws.Range("M132").Value = -18624.5
$ws.Range("N132").Value = -40149.00019999999

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1342.3158
$ws.Range("J135").Value = 2003.1428
$ws.Range("L135").Value = 18028.2852
$ws.Range("N135").Value = -23098.2852

$ws = $wb.Worksheets.Item("GSM")
# Row 17: Point of Honor
$ws.Range("H17").Value = 791.5714
$ws.Range("J17").Value = 1314.75
$ws.Range("L17").Value = 1314.75
$ws.Range("N17").Value = -1650.75

# Row 43: Get the Green Stuff
$ws.Range("H43").Value = 18666.533
$ws.Range("I43").Value = 13333
$ws.Range("J43").Value = 19999.916
$ws.Range("K43").Value = 13333
$ws.Range("L43").Value = 19999.916
$ws.Range("M43").Value = -13182
$ws.Range("N43").Value = -20301.916

# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 21200
$ws.Range("I46").Value = 30000
$ws.Range("J46").Value = 19000
$ws.Range("K46").Value = 30000
$ws.Range("L46").Value = 19000
$ws.Range("M46").Value = -29844
$ws.Range("N46").Value = -19312

# Row 99: Needle in a Hingan Stack
$ws.Range("H99").Value = 15152.143
$ws.Range("I99").Value = 5213
$ws.Range("K99").Value = 5213
$ws.Range("M99").Value = -2967

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 1962404.2
$ws.Range("I113").Value = 1210.6666
$ws.Range("K113").Value = 1210.6666
$ws.Range("M113").Value = 959.3334

# Row 132: On Board for Lar
$ws.Range("H132").Value = 4741
$ws.Range("J132").Value = 4302.875
$ws.Range("L132").Value = 12908.625
$ws.Range("N132").Value = -17968.625

$ws = $wb.Worksheets.Item("LTW")
# Row 75: Tally Ho, Chocobo
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872

# Row 78: Hunting Heretics (L)
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 1000
$ws.Range("J96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("N96").Value = -3746

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2384.5715
$ws.Range("I126").Value = 2038.4
$ws.Range("K126").Value = 6115.200000000001
$ws.Range("M126").Value = -3645.200000000001

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 680551
$ws.Range("I132").Value = 969.2549
$ws.Range("J132").Value = 3346602.5
$ws.Range("K132").Value = 2907.7647
$ws.Range("L132").Value = 10039807.5
$ws.Range("M132").Value = -377.7647000000002
$ws.Range("N132").Value = -10044867.5
